$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.982.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.260.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.58%  "
$ws.Range("E8").Value = "  +16.34%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0978"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.599.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.891"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.270.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.912.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("E24").Value = "  +4.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  +8.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0773"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.17%  "
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  +7.75%  "
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +1.59%  "
